# Daily attendance processing - 2025-11-13 20:48:38
# Normalises the "Recorded By" (column G) audit list on each attendance
# row: the most-recent recorder is tracked at the front of the list as
# edits land, so re-processing rotates the oldest entry (the original
# first recorder) to the back of the list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$col = 7  # Column G = "Recorded By"

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, $col)
    $current = $cell.Value2

    if ($current -ne $null -and $current -ne "") {
        $parts = $current -split ", "

        if ($parts.Count -gt 1) {
            $rotated = $parts[1..($parts.Count - 1)] + $parts[0]
            $newValue = $rotated -join ", "
            $cell.Value2 = $newValue
        }
    }
}
